# Update Mexico_M2 data: revise the last three existing monthly rows and
# append three new monthly rows (ECONOMICS:MXM2), mirroring open/high/low/close.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revise existing rows 450-452 (open/high/low/close columns C:F) ---
$revisions = @(
    @{ Row = 450; Value = 12492380715000 },
    @{ Row = 451; Value = 12552103031000 },
    @{ Row = 452; Value = 12748355850000 }
)

foreach ($rev in $revisions) {
    $r = $rev.Row
    $v = $rev.Value
    $ws.Cells.Item($r, 3).Value = $v   # C: open
    $ws.Cells.Item($r, 4).Value = $v   # D: high
    $ws.Cells.Item($r, 5).Value = $v   # E: low
    $ws.Cells.Item($r, 6).Value = $v   # F: close
}

# --- Append new rows 453-455 ---
$newRows = @(
    @{ Row = 453; Date = 45108.41666666666; Value = 12832232273000 },
    @{ Row = 454; Date = 45139.41666666666; Value = 12878456979000 },
    @{ Row = 455; Date = 45170.41666666666; Value = 12989512635000 }
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $prev = $r - 1

    # Copy formatting from the row above so the new datetime cell keeps the
    # same number format / font / border / alignment (style index) as the
    # rest of column A.
    $ws.Range("A$prev").Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    $ws.Cells.Item($r, 1).Value = $nr.Date            # A: datetime
    $ws.Cells.Item($r, 2).Value = "ECONOMICS:MXM2"    # B: symbol
    $ws.Cells.Item($r, 3).Value = $nr.Value           # C: open
    $ws.Cells.Item($r, 4).Value = $nr.Value           # D: high
    $ws.Cells.Item($r, 5).Value = $nr.Value           # E: low
    $ws.Cells.Item($r, 6).Value = $nr.Value           # F: close
    $ws.Cells.Item($r, 7).Value = 0                   # G: volume
}
